$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.5906356666666667
$ws.Cells.Item(2, 8).Value = 1.771907
$ws.Cells.Item(2, 9).Value = 0.294823169192623
$ws.Cells.Item(2, 10).Value = 0.294823169192623
$ws.Cells.Item(2, 13).Value = 82.43338033333333
$ws.Cells.Item(2, 14).Value = 247.300141
$ws.Cells.Item(2, 15).Value = 0.3670006993429558
$ws.Cells.Item(2, 16).Value = 0.3670006993429557
$ws.Cells.Item(2, 17).Value = 48.68809454876522
$ws.Cells.Item(2, 18).Value = 438.192850938887
$ws.Cells.Item(2, 19).Value = 0.1082003092761992
$ws.Cells.Item(2, 20).Value = 0.1082003092761992

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.5906356666666667
$ws.Cells.Item(3, 8).Value = 1.771907
$ws.Cells.Item(3, 9).Value = 0.294823169192623
$ws.Cells.Item(3, 10).Value = 0.294823169192623
$ws.Cells.Item(3, 15).Value = 0.3956886215996139
$ws.Cells.Item(3, 16).Value = 0.3956886215996139
$ws.Cells.Item(3, 17).Value = 52.49397359406522
$ws.Cells.Item(3, 18).Value = 472.4457623465869
$ws.Cells.Item(3, 19).Value = 0.1166581734334587
$ws.Cells.Item(3, 20).Value = 0.1166581734334587

# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.5906356666666667
$ws.Cells.Item(4, 8).Value = 1.771907
$ws.Cells.Item(4, 9).Value = 0.294823169192623
$ws.Cells.Item(4, 10).Value = 0.294823169192623
$ws.Cells.Item(4, 13).Value = 42.93483766666667
$ws.Cells.Item(4, 14).Value = 128.804513
$ws.Cells.Item(4, 15).Value = 0.1911496942879982
$ws.Cells.Item(4, 16).Value = 0.1911496942879981
$ws.Cells.Item(4, 17).Value = 25.35884646847678
$ws.Cells.Item(4, 18).Value = 228.229618216291
$ws.Cells.Item(4, 19).Value = 0.05635535866018864
$ws.Cells.Item(4, 20).Value = 0.05635535866018863

# Row 5
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.5906356666666667
$ws.Cells.Item(5, 8).Value = 1.771907
$ws.Cells.Item(5, 9).Value = 0.294823169192623
$ws.Cells.Item(5, 10).Value = 0.294823169192623
$ws.Cells.Item(5, 13).Value = 10.368389
$ws.Cells.Item(5, 14).Value = 31.105167
$ws.Cells.Item(5, 15).Value = 0.04616098476943217
$ws.Cells.Item(5, 16).Value = 0.04616098476943217
$ws.Cells.Item(5, 17).Value = 6.123940349274334
$ws.Cells.Item(5, 18).Value = 55.115463143469
$ws.Cells.Item(5, 19).Value = 0.01360932782277639
$ws.Cells.Item(5, 20).Value = 0.01360932782277639

# Row 6
$ws.Cells.Item(6, 9).Value = 0.2901829546991739
$ws.Cells.Item(6, 10).Value = 0.2901829546991739
$ws.Cells.Item(6, 13).Value = 82.43338033333333
$ws.Cells.Item(6, 14).Value = 247.300141
$ws.Cells.Item(6, 15).Value = 0.3670006993429558
$ws.Cells.Item(6, 16).Value = 0.3670006993429557
$ws.Cells.Item(6, 17).Value = 47.92179384518656
$ws.Cells.Item(6, 18).Value = 431.296144606679
$ws.Cells.Item(6, 19).Value = 0.1064973473120021
$ws.Cells.Item(6, 20).Value = 0.1064973473120021

# Row 7
$ws.Cells.Item(7, 9).Value = 0.2901829546991739
$ws.Cells.Item(7, 10).Value = 0.2901829546991739
$ws.Cells.Item(7, 15).Value = 0.3956886215996139
$ws.Cells.Item(7, 16).Value = 0.3956886215996139
$ws.Cells.Item(7, 19).Value = 0.1148220933566193
$ws.Cells.Item(7, 20).Value = 0.1148220933566193

# Row 8
$ws.Cells.Item(8, 9).Value = 0.2901829546991739
$ws.Cells.Item(8, 10).Value = 0.2901829546991739
$ws.Cells.Item(8, 13).Value = 42.93483766666667
$ws.Cells.Item(8, 14).Value = 128.804513
$ws.Cells.Item(8, 15).Value = 0.1911496942879982
$ws.Cells.Item(8, 16).Value = 0.1911496942879981
$ws.Cells.Item(8, 17).Value = 24.95972421752745
$ws.Cells.Item(8, 18).Value = 224.6375179577471
$ws.Cells.Item(8, 19).Value = 0.05546838307833511
$ws.Cells.Item(8, 20).Value = 0.05546838307833511

# Row 9
$ws.Cells.Item(9, 9).Value = 0.2901829546991739
$ws.Cells.Item(9, 10).Value = 0.2901829546991739
$ws.Cells.Item(9, 13).Value = 10.368389
$ws.Cells.Item(9, 14).Value = 31.105167
$ws.Cells.Item(9, 15).Value = 0.04616098476943217
$ws.Cells.Item(9, 16).Value = 0.04616098476943217
$ws.Cells.Item(9, 17).Value = 6.027555805130334
$ws.Cells.Item(9, 18).Value = 54.24800224617301
$ws.Cells.Item(9, 19).Value = 0.01339513095221739
$ws.Cells.Item(9, 20).Value = 0.01339513095221739

# Row 10
$ws.Cells.Item(10, 7).Value = 0.817256
$ws.Cells.Item(10, 8).Value = 2.451768
$ws.Cells.Item(10, 9).Value = 0.4079435387325965
$ws.Cells.Item(10, 10).Value = 0.4079435387325965
$ws.Cells.Item(10, 13).Value = 82.43338033333333
$ws.Cells.Item(10, 14).Value = 247.300141
$ws.Cells.Item(10, 15).Value = 0.3670006993429558
$ws.Cells.Item(10, 16).Value = 0.3670006993429557
$ws.Cells.Item(10, 17).Value = 67.36917467769867
$ws.Cells.Item(10, 18).Value = 606.322572099288
$ws.Cells.Item(10, 19).Value = 0.1497155640073031
$ws.Cells.Item(10, 20).Value = 0.149715564007303

# Row 11
$ws.Cells.Item(11, 7).Value = 0.817256
$ws.Cells.Item(11, 8).Value = 2.451768
$ws.Cells.Item(11, 9).Value = 0.4079435387325965
$ws.Cells.Item(11, 10).Value = 0.4079435387325965
$ws.Cells.Item(11, 15).Value = 0.3956886215996139
$ws.Cells.Item(11, 16).Value = 0.3956886215996139
$ws.Cells.Item(11, 17).Value = 72.63532716489865
$ws.Cells.Item(11, 18).Value = 653.717944484088
$ws.Cells.Item(11, 19).Value = 0.1614186165315698
$ws.Cells.Item(11, 20).Value = 0.1614186165315698

# Row 12
$ws.Cells.Item(12, 7).Value = 0.817256
$ws.Cells.Item(12, 8).Value = 2.451768
$ws.Cells.Item(12, 9).Value = 0.4079435387325965
$ws.Cells.Item(12, 10).Value = 0.4079435387325965
$ws.Cells.Item(12, 13).Value = 42.93483766666667
$ws.Cells.Item(12, 14).Value = 128.804513
$ws.Cells.Item(12, 15).Value = 0.1911496942879982
$ws.Cells.Item(12, 16).Value = 0.1911496942879981
$ws.Cells.Item(12, 17).Value = 35.08875369210934
$ws.Cells.Item(12, 18).Value = 315.798783228984
$ws.Cells.Item(12, 19).Value = 0.07797828271549996
$ws.Cells.Item(12, 20).Value = 0.07797828271549995

# Row 13
$ws.Cells.Item(13, 7).Value = 0.817256
$ws.Cells.Item(13, 8).Value = 2.451768
$ws.Cells.Item(13, 9).Value = 0.4079435387325965
$ws.Cells.Item(13, 10).Value = 0.4079435387325965
$ws.Cells.Item(13, 13).Value = 10.368389
$ws.Cells.Item(13, 14).Value = 31.105167
$ws.Cells.Item(13, 15).Value = 0.04616098476943217
$ws.Cells.Item(13, 16).Value = 0.04616098476943217
$ws.Cells.Item(13, 17).Value = 8.473628120583999
$ws.Cells.Item(13, 18).Value = 76.26265308525601
$ws.Cells.Item(13, 19).Value = 0.01883107547822365
$ws.Cells.Item(13, 20).Value = 0.01883107547822365

# Row 14
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.01412433333333333
$ws.Cells.Item(14, 8).Value = 0.042373
$ws.Cells.Item(14, 9).Value = 0.007050337375606629
$ws.Cells.Item(14, 10).Value = 0.007050337375606628
$ws.Cells.Item(14, 13).Value = 82.43338033333333
$ws.Cells.Item(14, 14).Value = 247.300141
$ws.Cells.Item(14, 15).Value = 0.3670006993429558
$ws.Cells.Item(14, 16).Value = 0.3670006993429557
$ws.Cells.Item(14, 17).Value = 1.164316541621444
$ws.Cells.Item(14, 18).Value = 10.478848874593
$ws.Cells.Item(14, 19).Value = 0.002587478747451412
$ws.Cells.Item(14, 20).Value = 0.002587478747451411

# Row 15
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.01412433333333333
$ws.Cells.Item(15, 8).Value = 0.042373
$ws.Cells.Item(15, 9).Value = 0.007050337375606629
$ws.Cells.Item(15, 10).Value = 0.007050337375606628
$ws.Cells.Item(15, 15).Value = 0.3956886215996139
$ws.Cells.Item(15, 16).Value = 0.3956886215996139
$ws.Cells.Item(15, 17).Value = 1.255329508321444
$ws.Cells.Item(15, 18).Value = 11.297965574893
$ws.Cells.Item(15, 19).Value = 0.002789738277966026
$ws.Cells.Item(15, 20).Value = 0.002789738277966026

# Row 16
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.01412433333333333
$ws.Cells.Item(16, 8).Value = 0.042373
$ws.Cells.Item(16, 9).Value = 0.007050337375606629
$ws.Cells.Item(16, 10).Value = 0.007050337375606628
$ws.Cells.Item(16, 13).Value = 42.93483766666667
$ws.Cells.Item(16, 14).Value = 128.804513
$ws.Cells.Item(16, 15).Value = 0.1911496942879982
$ws.Cells.Item(16, 16).Value = 0.1911496942879981
$ws.Cells.Item(16, 17).Value = 0.6064259588165557
$ws.Cells.Item(16, 18).Value = 5.457833629349
$ws.Cells.Item(16, 19).Value = 0.001347669833974454
$ws.Cells.Item(16, 20).Value = 0.001347669833974454

# Row 17
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.01412433333333333
$ws.Cells.Item(17, 8).Value = 0.042373
$ws.Cells.Item(17, 9).Value = 0.007050337375606629
$ws.Cells.Item(17, 10).Value = 0.007050337375606628
$ws.Cells.Item(17, 13).Value = 10.368389
$ws.Cells.Item(17, 14).Value = 31.105167
$ws.Cells.Item(17, 15).Value = 0.04616098476943217
$ws.Cells.Item(17, 16).Value = 0.04616098476943217
$ws.Cells.Item(17, 17).Value = 0.1464465823656667
$ws.Cells.Item(17, 18).Value = 1.318019241291
$ws.Cells.Item(17, 19).Value = 0.000325450516214736
$ws.Cells.Item(17, 20).Value = 0.0003254505162147359
